$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.0960744875850935
$ws.Range("E2").Value = 0.002520622277711838
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.02911706285278049
$ws.Range("E3").Value = 0.0007742253544701382
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 0.08259227143154145
$ws.Range("E4").Value = 0.001848109445836572
$ws.Range("B5").Value = 2470.366758096668
$ws.Range("D5").Value = 0.008296784368939976
$ws.Range("E5").Value = 0.000504062765808345
$ws.Range("B6").Value = 2481.984729487203
$ws.Range("D6").Value = 0.01369947817056663
$ws.Range("E6").Value = 0.000694477856152875
$ws.Range("B7").Value = 2497.913976174627
$ws.Range("D7").Value = 0.003373297067722395
$ws.Range("E7").Value = 0.0008187614242044651
$ws.Range("B8").Value = 2569.303537215177
$ws.Range("D8").Value = 0.01728622194449643
$ws.Range("E8").Value = 0.0006630908488656306
$ws.Range("B9").Value = 2593.831325411606
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 0.02208428679384985
$ws.Range("E9").Value = 0.001495566083159258
$ws.Range("B10").Value = 2609.227572242163
$ws.Range("D10").Value = 0.0268483955912421
$ws.Range("E10").Value = 0.001054856487612694
$ws.Range("B11").Value = 2621.783451658836
$ws.Range("D11").Value = 0.00309447226499785
$ws.Range("E11").Value = 0.0006268501198112199
$ws.Range("B12").Value = 2635.089616413751
$ws.Range("D12").Value = 0.01149130537220409
$ws.Range("E12").Value = 0.0006675264545959958
$ws.Range("B13").Value = 2649.075587546623
$ws.Range("D13").Value = 0.01472230434964875
$ws.Range("E13").Value = 0.0007573561146787284
$ws.Range("B14").Value = 2664.601302816965
$ws.Range("B15").Value = 2688.339468084317
$ws.Range("B16").Value = 2702.26533490114
$ws.Range("D16").Value = 0.004391319269637482
$ws.Range("E16").Value = 0.0007582133990740976
$ws.Range("B17").Value = 2715.640802344541
$ws.Range("B18").Value = 2737.363762576609
$ws.Range("D18").Value = 0.02780269134025979
$ws.Range("E18").Value = 0.00110310367550748
$ws.Range("B19").Value = 2749.680085527932
$ws.Range("B20").Value = 2767.781176406648
$ws.Range("D20").Value = 0.01853865151391864
$ws.Range("E20").Value = 0.0005810482725139638
$ws.Range("B21").Value = 2779.2
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0.07348957557846092
$ws.Range("E21").Value = 0.002774143209893858
$ws.Range("B22").Value = 2813.501259227846
$ws.Range("D22").Value = 0.003602772433087223
$ws.Range("E22").Value = 0.0002126147707377977
$ws.Range("B23").Value = 2828.359721907576
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 0.002374588601064437
$ws.Range("E23").Value = 0.000315379793147412
$ws.Range("B24").Value = 2853.114732511383
$ws.Range("D24").Value = 0.01742527498144794
$ws.Range("E24").Value = 0.0006691873879849242
$ws.Range("B25").Value = 2875.303509710963
$ws.Range("D25").Value = 0.01832457211464001
$ws.Range("E25").Value = 0.001638985953500033
$ws.Range("B26").Value = 2884.179211705924
$ws.Range("D26").Value = 0.01412043378719112
$ws.Range("E26").Value = 0.003734752414728085
$ws.Range("B27").Value = 2899.529675101117
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0.0189174637352376
$ws.Range("E27").Value = 0.002153207254417288
$ws.Range("B28").Value = 2917.634342911454
$ws.Range("D28").Value = 0.006523305033840098
$ws.Range("E28").Value = 0.001951855049495462
$ws.Range("B29").Value = 2928.524980912122
$ws.Range("D29").Value = 0.02974299713619296
$ws.Range("E29").Value = 0.001955421957050748
$ws.Range("B30").Value = 2949.317122042539
$ws.Range("D30").Value = 0.04751149746237109
$ws.Range("E30").Value = 0.001910788484899707
$ws.Range("B31").Value = 2970.664196645418
$ws.Range("D31").Value = 0.006987003767356034
$ws.Range("E31").Value = 0.0009924721260448912
$ws.Range("B32").Value = 2986.283965523739
$ws.Range("D32").Value = 0.01522743118189063
$ws.Range("E32").Value = 0.002702479254123934
$ws.Range("B33").Value = 2997.718657087334
$ws.Range("D33").Value = 0.01747133954907086
$ws.Range("E33").Value = 0.001217400296474966
$ws.Range("B34").Value = 3047.4
$ws.Range("D34").Value = 0.0331494547762802
$ws.Range("E34").Value = 0.002788763655782303
$ws.Range("B35").Value = 3064.813926466619
$ws.Range("D35").Value = 0.009817557518465999
$ws.Range("E35").Value = 0.001530694451803839
$ws.Range("B36").Value = 3078.112914411186
$ws.Range("D36").Value = 0.01021414327095403
$ws.Range("E36").Value = 0.0007919514663363961
$ws.Range("B37").Value = 3119.418862442509
$ws.Range("D37").Value = 0.01184523186921175
$ws.Range("E37").Value = 0.001547350109041175
$ws.Range("B38").Value = 3140.60165968082
$ws.Range("D38").Value = 0.003803189734194964
$ws.Range("E38").Value = 0.002303340261554697
$ws.Range("B39").Value = 3154.525784376794
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 0.001945554753831514
$ws.Range("E39").Value = 0.0004793395770309529
$ws.Range("B40").Value = 3171.655405406395
$ws.Range("D40").Value = 0.003706166898006592
$ws.Range("E40").Value = 0.0005053863951827171
$ws.Range("B41").Value = 3181.210075633211
$ws.Range("D41").Value = 0.02208980191773095
$ws.Range("E41").Value = 0.002659460181741564
$ws.Range("B42").Value = 3212.072026396334
$ws.Range("D42").Value = 0.002027262432905774
$ws.Range("E42").Value = 0.0004721022104027146
$ws.Range("B43").Value = 3231.273355870092
$ws.Range("D43").Value = 0.001160092827437504
$ws.Range("E43").Value = 0.0005248038981264897
$ws.Range("B44").Value = 3249.661277156506
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = 0.005866434030802336
$ws.Range("E44").Value = 0.0005087181939129394
$ws.Range("B45").Value = 3264.906563811741
$ws.Range("D45").Value = 0.006587785431038932
$ws.Range("E45").Value = 0.0004706482710866217
$ws.Range("B46").Value = 3286.170677130526
$ws.Range("D46").Value = 0.01082012791090631
$ws.Range("E46").Value = 0.0006869922483115117
$ws.Range("B47").Value = 3310.539257221813
$ws.Range("D47").Value = 0.003987701300459387
$ws.Range("E47").Value = 0.0003063996735747506
$ws.Range("B48").Value = 3327.874605893868
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0.01318027027306804
$ws.Range("E48").Value = 0.001289983899066233
$ws.Range("B49").Value = 3350.20085761773
$ws.Range("D49").Value = 0.01085366562236222
$ws.Range("E49").Value = 0.0004585738651058436
$ws.Range("B50").Value = 3368.24959104752
$ws.Range("D50").Value = 0.01357043037958504
$ws.Range("E50").Value = 0.0009480711635052563
$ws.Range("B51").Value = 3383.85849793409
$ws.Range("D51").Value = 0.006684249939476672
$ws.Range("E51").Value = 0.0005088623488223971
$ws.Range("B52").Value = 3401.382461109694
$ws.Range("D52").Value = 0.0130655937434392
$ws.Range("E52").Value = 0.0004628791821704674
$ws.Range("B53").Value = 3420.405841076768
$ws.Range("D53").Value = 0.008576425709640234
$ws.Range("E53").Value = 0.0003760143836651919
$ws.Range("B54").Value = 3441.040838634491
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 0.02324036883432738
$ws.Range("E54").Value = 0.00185922950674619
$ws.Range("B55").Value = 3466.321886486277
$ws.Range("D55").Value = 0.006071537232359859
$ws.Range("E55").Value = 0.001517884308089965
$ws.Range("B56").Value = 3474.201985506649
$ws.Range("C56").Value = 2
$ws.Range("D56").Value = 0.003528602170120224
$ws.Range("E56").Value = 0.0004670208754570885
$ws.Range("B57").Value = 3493.141363695906
$ws.Range("D57").Value = 0.02019293873546069
$ws.Range("E57").Value = 0.001429588583041465
$ws.Range("B58").Value = 3510.661101658533
$ws.Range("C58").Value = 4
$ws.Range("D58").Value = 0.008325312827123963
$ws.Range("E58").Value = 0.001800986040153347
$ws.Range("B59").Value = 3525.443535477975
$ws.Range("D59").Value = 0.005863072079783329
$ws.Range("E59").Value = 0.006102381144264281
$ws.Range("B60").Value = 3536.818978961481
$ws.Range("C60").Value = 1
$ws.Range("D60").Value = 0.01894605808449602
$ws.Range("E60").Value = 0.000777084289884611
$ws.Range("B61").Value = 3553.059086400061
$ws.Range("D61").Value = 0.009414025475432353
$ws.Range("E61").Value = 0.0008765386848633476
$ws.Range("B62").Value = 3572.945718819506
$ws.Range("D62").Value = 0.008934202554999708
$ws.Range("E62").Value = 0.0009266394039432185
$ws.Range("B63").Value = 3581.276241950112
$ws.Range("D63").Value = 0.03027168896261209
$ws.Range("E63").Value = 0.002136107084055289
